$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H4").Value = "Vocabulary code"
$ws.Range("H11").Value = "Vocabulary code"
$ws.Range("H19").Value = "Vocabulary code"

$ws.Range("E2").Value = "Generated code prefix"
$ws.Range("E9").Value = "Generated code prefix"
$ws.Range("E17").Value = "Generated code prefix"

$ws.Range("E17").Select()
